# Swap data between specific rows in the "Artfynd" sheet so that the
# Id / coordinates / time (and for some pairs species) columns line up
# with the committed state.
#
# Row pairs whose A,B,E,F,G,H,Q,R,Z,AB values are fully swapped:
#   8  <-> 9
#   11 <-> 12
# Row pairs where only A,Q,R,Z,AB are swapped (B,E,F,G,H already equal):
#   13 <-> 14
#   18 <-> 19

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

function Swap-Cell($ws, $col, $r1, $r2) {
    $addr1 = "$col$r1"
    $addr2 = "$col$r2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

# Full swap columns for row pairs 8/9 and 11/12
$fullCols = @("A","B","E","F","G","H","Q","R","Z","AB")
foreach ($col in $fullCols) {
    Swap-Cell $ws $col 8 9
    Swap-Cell $ws $col 11 12
}

# Partial swap columns (only identity/location/time) for row pairs 13/14 and 18/19
$partialCols = @("A","Q","R","Z","AB")
foreach ($col in $partialCols) {
    Swap-Cell $ws $col 13 14
    Swap-Cell $ws $col 18 19
}
